$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "release/8.0.2"
$ws.Range("B4").Value = "X"
$ws.Range("C4").Value = "X"
$ws.Range("D4").Value = "X"
